$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1422776694794123
$ws.Range("D2").Value = 0.1016491568980555
$ws.Range("E2").Value = 0.1314834985687696
$ws.Range("F2").Value = 2.228370127467002
$ws.Range("G2").Value = 1.558939120223002
$ws.Range("H2").Value = 1.393057684661827
$ws.Range("J2").Value = 0.1868716232151826
$ws.Range("K2").Value = 0.7288204472436064
$ws.Range("M2").Value = 0.2998801559184585

$ws.Range("B3").Value = 0.1330581801436637
$ws.Range("D3").Value = 0.1006582097491844
$ws.Range("E3").Value = 0.1310989901654978
$ws.Range("F3").Value = 2.22221473876165
$ws.Range("G3").Value = 1.55099031119093
$ws.Range("H3").Value = 1.395134461948487
$ws.Range("J3").Value = 0.1869436934763442
$ws.Range("K3").Value = 0.6597646867518847
$ws.Range("M3").Value = 0.2846789530053826

$ws.Range("B4").Value = 0.1274702943850912
$ws.Range("D4").Value = 0.1000821577253319
$ws.Range("E4").Value = 0.1309151267613977
$ws.Range("F4").Value = 2.219604470070266
$ws.Range("G4").Value = 1.547038132775782
$ws.Range("H4").Value = 1.397051250185015
$ws.Range("J4").Value = 0.1870697789643891
$ws.Range("K4").Value = 0.617574962434702
$ws.Range("M4").Value = 0.2754861735595426

$ws.Range("B5").Value = 0.1252116289615799
$ws.Range("D5").Value = 0.09985560076377453
$ws.Range("E5").Value = 0.1308533523143609
$ws.Range("F5").Value = 2.218834446304101
$ws.Range("G5").Value = 1.545660552036907
$ws.Range("H5").Value = 1.397993605551193
$ws.Range("J5").Value = 0.1871417353583773
$ws.Range("K5").Value = 0.6004356922380225
$ws.Range("M5").Value = 0.2717755601349978

$ws.Range("B6").Value = 0.1248376964037163
$ws.Range("D6").Value = 0.09981847702515267
$ws.Range("E6").Value = 0.1308438895792534
$ws.Range("F6").Value = 2.218724312983284
$ws.Range("G6").Value = 1.545445861799152
$ws.Range("H6").Value = 1.398159819165457
$ws.Range("J6").Value = 0.1871549263128038
$ws.Range("K6").Value = 0.5975929652956324
$ws.Range("M6").Value = 0.271161564115225

$ws.Range("B7").Value = 0.127439758413729
$ws.Range("D7").Value = 0.1000790690838897
$ws.Range("E7").Value = 0.1309142403749775
$ws.Range("F7").Value = 2.219592896578206
$ws.Range("G7").Value = 1.547018611640183
$ws.Range("H7").Value = 1.397063306365254
$ws.Range("J7").Value = 0.1870706660865658
$ws.Range("K7").Value = 0.6173435995756904
$ws.Range("M7").Value = 0.2754359869689154

$ws.Range("B8").Value = 0.1390837166955379
$ws.Range("D8").Value = 0.1013007804004644
$ws.Range("E8").Value = 0.1313400912380125
$ws.Range("F8").Value = 2.226004892358503
$ws.Range("G8").Value = 1.556005372144995
$ws.Range("H8").Value = 1.393640524311863
$ws.Range("J8").Value = 0.1868794865473831
$ws.Range("K8").Value = 0.7049664731369489
$ws.Range("M8").Value = 0.2946096126244342

$ws.Range("B9").Value = 0.1624926424440076
$ws.Range("D9").Value = 0.1039518931783974
$ws.Range("E9").Value = 0.1325889382863217
$ws.Range("F9").Value = 2.247873270061987
$ws.Range("G9").Value = 1.58102027915163
$ws.Range("H9").Value = 1.392025249998994
$ws.Range("J9").Value = 0.1871542048314794
$ws.Range("K9").Value = 0.8784602220005979
$ws.Range("M9").Value = 0.3333239896289228

$ws.Range("B10").Value = 0.1800391661567886
$ws.Range("D10").Value = 0.1060533761812295
$ws.Range("E10").Value = 0.1337581149180416
$ws.Range("F10").Value = 2.269635669382993
$ws.Range("G10").Value = 1.603944273634738
$ws.Range("H10").Value = 1.393955417394238
$ws.Range("J10").Value = 0.1877527782347315
$ws.Range("K10").Value = 1.00694642831786
$ws.Range("M10").Value = 0.3624474370489139

$ws.Range("B11").Value = 0.1880966577398482
$ws.Range("D11").Value = 0.1070423852237425
$ws.Range("E11").Value = 0.1343445438315989
$ws.Range("F11").Value = 2.280779437789647
$ws.Range("G11").Value = 1.615368743724957
$ws.Range("H11").Value = 1.395512566531039
$ws.Range("J11").Value = 0.1881114087213831
$ws.Range("K11").Value = 1.065621847957971
$ws.Range("M11").Value = 0.3758444316270868

$ws.Range("B12").Value = 0.1911585952388606
$ws.Range("D12").Value = 0.107421609821543
$ws.Range("E12").Value = 0.1345744421928501
$ws.Range("F12").Value = 2.285178622957346
$ws.Range("G12").Value = 1.619838775055058
$ws.Range("H12").Value = 1.396200027225063
$ws.Range("J12").Value = 0.188259637426043
$ws.Range("K12").Value = 1.087873144029061
$ws.Range("M12").Value = 0.3809388489171184

$ws.Range("B13").Value = 0.1904986763660759
$ws.Range("D13").Value = 0.1073397282767061
$ws.Range("E13").Value = 0.1345245814559313
$ws.Range("F13").Value = 2.284223199642966
$ws.Range("G13").Value = 1.618869667657322
$ws.Range("H13").Value = 1.396047618122793
$ws.Range("J13").Value = 0.1882271610704151
$ws.Range("K13").Value = 1.083079503441866
$ws.Range("M13").Value = 0.379840730762858

$ws.Range("B14").Value = 0.1883483505876882
$ws.Range("D14").Value = 0.1070734901359245
$ws.Range("E14").Value = 0.134363300864198
$ws.Range("F14").Value = 2.281137766000711
$ws.Range("G14").Value = 1.615733609973347
$ws.Range("H14").Value = 1.39556716338663
$ws.Range("J14").Value = 0.1881233545615686
$ws.Range("K14").Value = 1.067451832442316
$ws.Range("M14").Value = 0.3762631267664034

$ws.Range("B15").Value = 0.187032609251915
$ws.Range("D15").Value = 0.106911023516723
$ws.Range("E15").Value = 0.1342655311347229
$ws.Range("F15").Value = 2.279271208640779
$ws.Range("G15").Value = 1.61383143243205
$ws.Range("H15").Value = 1.395285612110939
$ws.Range("J15").Value = 0.1880613881769264
$ws.Range("K15").Value = 1.057883613851374
$ws.Range("M15").Value = 0.3740745075512066

$ws.Range("B16").Value = 0.179514099343919
$ws.Range("D16").Value = 0.1059894032487563
$ws.Range("E16").Value = 0.1337208870274402
$ws.Range("F16").Value = 2.268932463749707
$ws.Range("G16").Value = 1.60321775353114
$ws.Range("H16").Value = 1.393867331960507
$ws.Range("J16").Value = 0.1877310788827558
$ws.Range("K16").Value = 1.003116375651274
$ws.Range("M16").Value = 0.3615748932852512

$ws.Range("B17").Value = 0.1749209899058854
$ws.Range("D17").Value = 0.1054324498701718
$ws.Range("E17").Value = 0.1334007294537578
$ws.Range("F17").Value = 2.262908869469214
$ws.Range("G17").Value = 1.596962172999724
$ws.Range("H17").Value = 1.393171302890693
$ws.Range("J17").Value = 0.1875505633296584
$ws.Range("K17").Value = 0.9695761695142266
$ws.Range("M17").Value = 0.3539447699056311

$ws.Range("B18").Value = 0.1722862671797287
$ws.Range("D18").Value = 0.1051152167975786
$ws.Range("E18").Value = 0.1332217200992964
$ws.Range("F18").Value = 2.259561318897724
$ws.Range("G18").Value = 1.593457863898607
$ws.Range("H18").Value = 1.392834877142832
$ws.Range("J18").Value = 0.1874548614605658
$ws.Range("K18").Value = 0.9503060610770717
$ws.Range("M18").Value = 0.3495701168195922

$ws.Range("B19").Value = 0.1713954203941057
$ws.Range("D19").Value = 0.1050083428649629
$ws.Range("E19").Value = 0.1331619933471799
$ws.Range("F19").Value = 2.258447988954558
$ws.Range("G19").Value = 1.592287446674618
$ws.Range("H19").Value = 1.392731941614272
$ws.Range("J19").Value = 0.1874238538413593
$ws.Range("K19").Value = 0.943785208511116
$ws.Range("M19").Value = 0.3480913407349604

$ws.Range("B20").Value = 0.175409198980546
$ws.Range("D20").Value = 0.1054914167550791
$ws.Range("E20").Value = 0.1334342792490339
$ws.Range("F20").Value = 2.263537973002641
$ws.Range("G20").Value = 1.597618384794657
$ws.Range("H20").Value = 1.39323878067421
$ws.Range("J20").Value = 0.1875689384605082
$ws.Range("K20").Value = 0.9731443767427663
$ws.Range("M20").Value = 0.3547555621810687

$ws.Range("B21").Value = 0.1889796628205005
$ws.Range("D21").Value = 0.1071515632796931
$ws.Range("E21").Value = 0.1344104604923899
$ws.Range("F21").Value = 2.282039163759379
$ws.Range("G21").Value = 1.6166508379672
$ws.Range("H21").Value = 1.39570562908446
$ws.Range("J21").Value = 0.1881535078246017
$ws.Range("K21").Value = 1.072041187528953
$ws.Range("M21").Value = 0.3773133794159591

$ws.Range("B22").Value = 0.197911294872938
$ws.Range("D22").Value = 0.1082639897035094
$ws.Range("E22").Value = 0.1350940860689249
$ws.Range("F22").Value = 2.29517594300853
$ws.Range("G22").Value = 1.62992828015274
$ws.Range("H22").Value = 1.397887992218358
$ws.Range("J22").Value = 0.1886079730349834
$ws.Range("K22").Value = 1.136863425496756
$ws.Range("M22").Value = 0.3921801278104837

$ws.Range("B23").Value = 0.1931386294456274
$ws.Range("D23").Value = 0.107667771137578
$ws.Range("E23").Value = 0.1347250518379219
$ws.Range("F23").Value = 2.288068828319865
$ws.Range("G23").Value = 1.622764933730281
$ws.Range("H23").Value = 1.396671007430768
$ws.Range("J23").Value = 0.1883587876817074
$ws.Range("K23").Value = 1.102249543233086
$ws.Range("M23").Value = 0.3842341592184582

$ws.Range("B24").Value = 0.175188461219193
$ws.Range("D24").Value = 0.1054647485833513
$ws.Range("E24").Value = 0.1334190956448289
$ws.Range("F24").Value = 2.263253195554682
$ws.Range("G24").Value = 1.597321424567326
$ws.Range("H24").Value = 1.393208075448939
$ws.Range("J24").Value = 0.1875606058995345
$ws.Range("K24").Value = 0.9715311509751245
$ws.Range("M24").Value = 0.3543889655879795

$ws.Range("B25").Value = 0.1560986075424751
$ws.Range("D25").Value = 0.1032075646217763
$ws.Range("E25").Value = 0.1322068583109477
$ws.Range("F25").Value = 2.240959211369599
$ws.Range("G25").Value = 1.573457678323507
$ws.Range("H25").Value = 1.391915573598965
$ws.Range("J25").Value = 0.1870102670007228
$ws.Range("K25").Value = 0.8313467902893024
$ws.Range("M25").Value = 0.3227313426600773
